# Auto-generated Excel COM-interop edit script
# Applies numeric updates to the FFXIV leve-profit tables across 7 crafter sheets
# (ALC, ARM, BSM, CRP, GSM, LTW, WVR) as captured by the upstream data-refresh diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 655.08826
$ws.Range("I15").Value = 655.08826
$ws.Range("K15").Value = 1965.26478
$ws.Range("M15").Value = -1796.26478

# Row 125
$ws.Range("H125").Value = 1176.2
$ws.Range("J125").Value = 1160.3334
$ws.Range("L125").Value = 10443.0006
$ws.Range("N125").Value = -15363.0006

# Row 137
$ws.Range("H137").Value = 13890530
$ws.Range("I137").Value = 47620376
$ws.Range("J137").Value = 1769.4706
$ws.Range("K137").Value = 142861128
$ws.Range("L137").Value = 5308.4118
$ws.Range("M137").Value = -142858578
$ws.Range("N137").Value = -10408.4118

# Row 138
$ws.Range("H138").Value = 3797.577
$ws.Range("I138").Value = 5713.6665
$ws.Range("J138").Value = 3222.75
$ws.Range("K138").Value = 17140.9995
$ws.Range("L138").Value = 9668.25
$ws.Range("M138").Value = -12000.9995
$ws.Range("N138").Value = -19948.25

# Row 141
$ws.Range("H141").Value = 3631.2
$ws.Range("I141").Value = 3556.7
$ws.Range("K141").Value = 10670.1
$ws.Range("M141").Value = -5490.099999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1114.1818
$ws.Range("I2").Value = 648.86664
$ws.Range("K2").Value = 648.86664
$ws.Range("M2").Value = -535.86664

# Row 32
$ws.Range("H32").Value = 13207.214
$ws.Range("I32").Value = 9490.5
$ws.Range("K32").Value = 9490.5
$ws.Range("M32").Value = -9203.5

# Row 40
$ws.Range("H40").Value = 30000
$ws.Range("J40").Value = 30000
$ws.Range("L40").Value = 30000
$ws.Range("N40").Value = -30352

# Row 74
$ws.Range("H74").Value = 2116.25
$ws.Range("I74").Value = 1541.05
$ws.Range("J74").Value = 4992.25
$ws.Range("K74").Value = 1541.05
$ws.Range("L74").Value = 4992.25
$ws.Range("M74").Value = -667.05
$ws.Range("N74").Value = -6740.25

# Row 77
$ws.Range("H77").Value = 2116.25
$ws.Range("I77").Value = 1541.05
$ws.Range("J77").Value = 4992.25
$ws.Range("K77").Value = 7705.25
$ws.Range("L77").Value = 24961.25
$ws.Range("M77").Value = -3337.25
$ws.Range("N77").Value = -33697.25

# Row 86
$ws.Range("H86").Value = 32000
$ws.Range("I86").Value = 32000
$ws.Range("K86").Value = 32000
$ws.Range("M86").Value = -30814

# Row 89
$ws.Range("H89").Value = 32000
$ws.Range("I89").Value = 32000
$ws.Range("K89").Value = 96000
$ws.Range("M89").Value = -90072

# Row 105
$ws.Range("H105").Value = 33250
$ws.Range("J105").Value = 33250
$ws.Range("L105").Value = 33250
$ws.Range("N105").Value = -40238

# Row 116
$ws.Range("H116").Value = 1114.1818
$ws.Range("I116").Value = 648.86664
$ws.Range("K116").Value = 648.86664
$ws.Range("M116").Value = 1645.13336

# Row 122
$ws.Range("H122").Value = 2246.1785
$ws.Range("I122").Value = 2335.3076
$ws.Range("K122").Value = 7005.9228
$ws.Range("M122").Value = -4555.9228

# Row 124
$ws.Range("H124").Value = 99369.28999999999
$ws.Range("J124").Value = 99369.28999999999
$ws.Range("L124").Value = 99369.28999999999
$ws.Range("N124").Value = -109189.29

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1114.1818
$ws.Range("I3").Value = 648.86664
$ws.Range("K3").Value = 648.86664
$ws.Range("M3").Value = -534.86664

# Row 105
$ws.Range("H105").Value = 2341.8572
$ws.Range("I105").Value = 2315.5
$ws.Range("K105").Value = 2315.5
$ws.Range("M105").Value = -568.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1345.6364
$ws.Range("I31").Value = 1456.4445
$ws.Range("J31").Value = 847
$ws.Range("K31").Value = 1456.4445
$ws.Range("L31").Value = 847
$ws.Range("M31").Value = -1161.4445
$ws.Range("N31").Value = -1437

# Row 34
$ws.Range("H34").Value = 1345.6364
$ws.Range("I34").Value = 1456.4445
$ws.Range("J34").Value = 847
$ws.Range("K34").Value = 1456.4445
$ws.Range("L34").Value = 847
$ws.Range("M34").Value = -1254.4445
$ws.Range("N34").Value = -1251

# Row 86
$ws.Range("H86").Value = 12499.538
$ws.Range("I86").Value = 10856.111
$ws.Range("K86").Value = 10856.111
$ws.Range("M86").Value = -9733.111000000001

# Row 89
$ws.Range("H89").Value = 12499.538
$ws.Range("I89").Value = 10856.111
$ws.Range("K89").Value = 54280.55500000001
$ws.Range("M89").Value = -48664.55500000001

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 59999
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 59999
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 59999
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -60345

# Row 24
$ws.Range("H24").Value = 11849.75
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 11849.75
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 11849.75
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -12195.75

# Row 30
$ws.Range("H30").Value = 59999
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 59999
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 59999
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -60209

# Row 95
$ws.Range("H95").Value = 9666.333000000001
$ws.Range("J95").Value = 9499.5
$ws.Range("L95").Value = 9499.5
$ws.Range("N95").Value = -14991.5

# Row 98
$ws.Range("H98").Value = 10574.75
$ws.Range("J98").Value = 11433
$ws.Range("L98").Value = 11433
$ws.Range("N98").Value = -17423

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 136
$ws.Range("H136").Value = 64765.934
$ws.Range("J136").Value = 64765.934
$ws.Range("L136").Value = 194297.802
$ws.Range("N136").Value = -199397.802

$ws = $wb.Worksheets.Item("LTW")
# Row 38
$ws.Range("H38").Value = 79475
$ws.Range("J38").Value = 79475
$ws.Range("L38").Value = 79475
$ws.Range("N38").Value = -80295

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

# Row 105
$ws.Range("H105").Value = 34999.332
$ws.Range("J105").Value = 34999.332
$ws.Range("L105").Value = 34999.332
$ws.Range("N105").Value = -41987.332

# Row 127
$ws.Range("H127").Value = 87500
$ws.Range("J127").Value = 87500
$ws.Range("L127").Value = 87500
$ws.Range("N127").Value = -97420

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 14999
$ws.Range("J18").Value = 14999
$ws.Range("L18").Value = 14999
$ws.Range("N18").Value = -15345

# Row 42
$ws.Range("H42").Value = 30000
$ws.Range("J42").Value = 30000
$ws.Range("L42").Value = 30000
$ws.Range("N42").Value = -30756

# Row 64
$ws.Range("H64").Value = 90000
$ws.Range("J64").Value = 90000
$ws.Range("L64").Value = 90000
$ws.Range("N64").Value = -90496

# Row 67
$ws.Range("H67").Value = 90000
$ws.Range("J67").Value = 90000
$ws.Range("L67").Value = 90000
$ws.Range("N67").Value = -91716

# Row 94
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51802
